# "cierre 23 julio 2022" - weekly payroll receipt roll-forward.
# Update the week label and the hours total for the new week; every other
# touched cell (H9, B27, H27, B43, dates, SUM) is a formula that recalculates
# on its own once these source values change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New week label (drives B9, and transitively H9/B27/H27/B43 via formulas).
$ws.Range("B9").Value = "SEMANA   29  DEL    18      Al   24   DE   JULIO          2022"

# Hours worked entered for the new week.
$ws.Range("K21").Value = 560

# Last user action was on the hours cell below it.
$ws.Range("K22").Select() | Out-Null
